$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: out_vars  -> add row 18 (2020-06-17 raw national data)
# ---------------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")
$wsOut.Range("A17:J17").Copy()
$wsOut.Range("A18:J18").PasteSpecial(-4122)

$outVals = @(43999,159793,222801,59076,19080,32.055221442741548,51222,4654,4828,441670)
for ($c = 1; $c -le 10; $c++) {
  $wsOut.Cells.Item(18, $c).Value = $outVals[$c-1]
}

# ---------------------------------------------------------------------------
# Sheet: dates_dx -> fill row 18 (already present, blank)
# ---------------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")
$dxVals = @(43999,0,1,1,1,1,0,0,1,0,4)
for ($c = 1; $c -le 11; $c++) {
  $wsDx.Cells.Item(18, $c).Value = $dxVals[$c-1]
}
$wsDx.Range("C29").Select()

# ---------------------------------------------------------------------------
# Sheet: dates_sx -> add row 18
# ---------------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")
$wsSx.Range("A17:M17").Copy()
$wsSx.Range("A18:M18").PasteSpecial(-4122)
$sxVals = @(43999,0,1,1,0,1,1,1,0,1,1,0,0)
for ($c = 1; $c -le 13; $c++) {
  $wsSx.Cells.Item(18, $c).Value = $sxVals[$c-1]
}
$wsSx.Range("D20").Select()

# ---------------------------------------------------------------------------
# Sheet: dates_deaths -> fill row 18 (A18 already exists but w/ wrong style)
# ---------------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")
$wsDeaths.Cells.Item(17,1).Copy()
$wsDeaths.Cells.Item(18,1).PasteSpecial(-4122)
$deathsVals = @(43999,0,0,0,0,2,1,1,1,2)
for ($c = 1; $c -le 10; $c++) {
  $wsDeaths.Cells.Item(18, $c).Value = $deathsVals[$c-1]
}
$wsDeaths.Range("I25").Select()

# ---------------------------------------------------------------------------
# Sheet: control_obs -> add column R (2020-06-17) + extend S:W blank columns
# ---------------------------------------------------------------------------
$wsCtl = $wb.Worksheets.Item("control_obs")

$wsCtl.Range("Q1:Q20").Copy()
$wsCtl.Range("R1:R20").PasteSpecial(-4122)
$wsCtl.Range("R1:R20").Copy()
$wsCtl.Range("S1:W20").PasteSpecial(-4122)

$rRows    = @(1,2,3,4,5,6,7,8,10,11,12,13,14,15,16,18)
$rValues  = @(43999,3561,3376,3376,3376,3376,2567,5198,157,157,157,157,157,92,169,828)
for ($i = 0; $i -lt $rRows.Length; $i++) {
  $wsCtl.Cells.Item($rRows[$i], 18).Value = $rValues[$i]
}

$wsCtl.Range("C20:R20").FormulaR1C1 = "=SUM(R[-18]C:R[-2]C)"

$wsCtl.Range("R25").Select()

# ---------------------------------------------------------------------------
# Sheet: control_obs_mpio -> selection only
# ---------------------------------------------------------------------------
$wsMpio = $wb.Worksheets.Item("control_obs_mpio")
$wsMpio.Range("B2").Select()

# ---------------------------------------------------------------------------
# Sheet: anomalias -> add notes about June 17th deaths-series anomaly
# ---------------------------------------------------------------------------
$wsAnom = $wb.Worksheets.Item("anomalias")

$wsAnom.Cells.Item(8,1).Value = "17 de junio"
$wsAnom.Range("B8:B11").Merge()
$wsAnom.Range("B8:B11").WrapText = $true
$wsAnom.Range("B8:B11").VerticalAlignment = -4108
$wsAnom.Cells.Item(8,2).Value = "En serie de defunciones el día 12 de junio de 2020 la serie comenzó el 5 de febrero de 2020, para el 13 de junio la serie comenzó el 6 de marzo por lo cual hay una diferencia de 29 días, es decir 29 observaciones"

$wsAnom.Cells.Item(12,1).Value = "17 de junio"
$wsAnom.Range("B12:B14").Merge()
$wsAnom.Range("B12:B14").WrapText = $true
$wsAnom.Range("B12:B14").HorizontalAlignment = -4108
$wsAnom.Range("B12:B14").VerticalAlignment = -4108
$wsAnom.Cells.Item(12,2).Value = "En serie de defunciones el día 14 de junio, la serie comenzó el 6 de marzo, el día 15 de junio comenzó el 18 de marzo; lo cual es una diferencia de 11 días es decir 11 observaciones."

$wsAnom.Range("B15").Select()

# ---------------------------------------------------------------------------
# Final active sheet/selection must be out_vars, A18 (matches target workbook
# bookViews activeTab omission == 0 / first sheet, tabSelected on sheet1).
# ---------------------------------------------------------------------------
$wsOut.Activate()
$wsOut.Range("A18").Select()
